$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'45.884.63"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").Value = "'2.373.75"
$ws.Range("E3").Value = "  +2.85%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'300.01"
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").Value = "'98.30"
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("D7").Value = "'0.564"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = "  -4.48%  "
$ws.Range("D10").Value = "'34.31"
$ws.Range("E10").Value = "  -7.12%  "
$ws.Range("D11").Value = "'0.0785"
$ws.Range("E11").Value = "  -2.22%  "
$ws.Range("D12").Value = "'7.10"
$ws.Range("E12").Value = "  -4.58%  "
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "'2.744.75"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").Value = "'2.378.52"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "'0.822"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "'13.71"
$ws.Range("E17").Value = "  -2.29%  "
$ws.Range("D18").Value = "'45.835.18"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("E19").Value = "  -6.86%  "
$ws.Range("D20").Value = "'0.0₃0944"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("D21").Value = "'6.01"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("D22").Value = "'66.77"
$ws.Range("E22").Value = "  -0.17%  "
$ws.Range("D23").Value = "'242.92"
$ws.Range("E23").Value = "  -2.68%  "
$ws.Range("D24").Value = "'2.78"
$ws.Range("E24").Value = "  -5.50%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").Value = "'38.91"
$ws.Range("E27").Value = "  -12.00%  "
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("D29").Value = "'9.70"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("D30").Value = "'20.98"
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("D31").Value = "'3.71"
$ws.Range("E31").Value = "  +15.79%  "
$ws.Range("D32").Value = "'2.81"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").Value = "'5.52"
$ws.Range("E33").Value = "  -4.83%  "
$ws.Range("D34").Value = "'147.45"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "'0.0767"
$ws.Range("E35").Value = "  -5.27%  "
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Value = "'1.91"
$ws.Range("E37").Value = "  +5.07%  "
$ws.Range("D38").Value = "'0.115"
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("D39").Value = "'14.88"
$ws.Range("E39").Value = "  -8.23%  "
$ws.Range("D40").Value = "'3.85"
$ws.Range("E40").Value = "  -4.49%  "
$ws.Range("E41").Value = "  -2.96%  "
$ws.Range("D42").Value = "'3.19"
$ws.Range("E42").Value = "  -8.01%  "
$ws.Range("D43").Value = "'1.943.65"
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.09%  "
$ws.Range("D45").Value = "'94.53"
$ws.Range("E45").Value = "  +6.41%  "
$ws.Range("D46").Value = "'1.80"
$ws.Range("E46").Value = "  -10.72%  "
$ws.Range("D47").Value = "'8.47"
$ws.Range("E47").Value = "  +5.05%  "
$ws.Range("E48").Value = "  -6.09%  "
$ws.Range("D49").Value = "'98.85"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").Value = "'2.608.29"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("D51").Value = "'68.39"
